$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 124, shifting existing rows 124:224 down to 125:225
$ws.Rows(124).Insert()

# Populate the newly inserted row 124 with the new record
$ws.Cells.Item(124,1).Value  = 10
$ws.Cells.Item(124,2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(124,3).Value  = 'La Araucanía'
$ws.Cells.Item(124,4).Value  = 44447
$ws.Cells.Item(124,5).Value  = 9
$ws.Cells.Item(124,6).Value  = 100112032
$ws.Cells.Item(124,7).Value  = 'Zapallo italiano'
$ws.Cells.Item(124,8).Value  = 'Sin especificar'
$ws.Cells.Item(124,9).Value  = 'Primera'
$ws.Cells.Item(124,10).Value = 95
$ws.Cells.Item(124,11).Value = 16000
$ws.Cells.Item(124,12).Value = 16000
$ws.Cells.Item(124,13).Value = 16000
$ws.Cells.Item(124,14).Value = '$/caja 60 unidades'
$ws.Cells.Item(124,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(124,16).Value = 267
$ws.Cells.Item(124,17).Value = 60
$ws.Cells.Item(124,18).Value = 'Hortaliza'
